# Update gh-pages output (456a3b4) — refresh scraped counters / swap in new
# rows for 上海-漫展信息.xlsx across all four sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (1)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 33
$ws.Range("F5").Value = 142
$ws.Range("F7").Value = 289
$ws.Range("F8").Value = 352
$ws.Range("F9").Value = 3330
$ws.Range("F10").Value = 1174
$ws.Range("F11").Value = 1053
$ws.Range("F13").Value = 97
$ws.Range("F15").Value = 1519
$ws.Range("F17").Value = 837
$ws.Range("F18").Value = 1720
$ws.Range("I18").Value = "//i1.hdslb.com/bfs/openplatform/202410/GwoLrudT1728526473661.jpeg"
$ws.Range("F20").Value = 395
$ws.Range("F24").Value = 115

# ---------------------------------------------------------------------
# Sheet "演出" (2)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")

# Row 36 now shows the event that used to be row 37's
# (leading apostrophe keeps the date-shaped string as text, matching the
# workbook's existing plain-text date column instead of letting Excel coerce
# it into a date serial)
$ws.Range("B36").Value = "'2024-11-15"
$ws.Range("C36").Value = "上海·“法国姐姐”乔伊丝·乔纳森《小意思》"
$ws.Range("D36").Value = "高青西路777号 上海前滩31演艺中心"
$ws.Range("E36").Value = "2024.11.15 19:30-11.15 21:00"
$ws.Range("F36").Value = 4
$ws.Range("G36").Value = 280
$ws.Range("H36").Value = "https://show.bilibili.com/platform/detail.html?id=91619"
$ws.Range("I36").Value = "//i1.hdslb.com/bfs/openplatform/202408/VnZEk71H1725014748758.jpeg"

# Row 37 now shows the event that used to be row 38's
$ws.Range("B37").Value = "'2024-11-16"
$ws.Range("C37").Value = "上海·变形金刚音乐会40周年特别版"
$ws.Range("E37").Value = "2024.11.16 19:30-11.16 21:30"
$ws.Range("F37").Value = 54
$ws.Range("G37").Value = 266
$ws.Range("H37").Value = "https://show.bilibili.com/platform/detail.html?id=90031"
$ws.Range("I37").Value = "//i1.hdslb.com/bfs/openplatform/202409/5zTUqO9f1727061199503.jpeg"

# Row 38 becomes a brand-new event
$ws.Range("C38").Value = "上海·钢琴&大提琴烛光音乐会演奏贝加尔湖畔&权力的游戏&久石让曲目经典演奏"
$ws.Range("D38").Value = "南苏州路1247号 八号桥艺术空间"
$ws.Range("E38").Value = "2024.11.16 13:20-12.24 21:40"
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 198
$ws.Range("H38").Value = "https://show.bilibili.com/platform/detail.html?id=93222"
$ws.Range("I38").Value = "//i0.hdslb.com/bfs/openplatform/202410/0l20Sp1l1728459887251.jpeg"

$ws.Range("F40").Value = 368
$ws.Range("F48").Value = 302

# ---------------------------------------------------------------------
# Sheet "本地生活" (3)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 2512
$ws.Range("F7").Value = 9551
$ws.Range("F12").Value = 2836
$ws.Range("F13").Value = 385
$ws.Range("F14").Value = 702

# ---------------------------------------------------------------------
# Sheet "全部类型" (4)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 2836
$ws.Range("F7").Value = 385
$ws.Range("F8").Value = 142
$ws.Range("F9").Value = 702
$ws.Range("F15").Value = 289
$ws.Range("F16").Value = 352
$ws.Range("F17").Value = 1174
$ws.Range("F19").Value = 1053
$ws.Range("F20").Value = 97
$ws.Range("F23").Value = 1519
$ws.Range("F26").Value = 837
$ws.Range("F29").Value = 1720
$ws.Range("I29").Value = "//i1.hdslb.com/bfs/openplatform/202410/GwoLrudT1728526473661.jpeg"
$ws.Range("F30").Value = 395
$ws.Range("F42").Value = 122
$ws.Range("F43").Value = 368
$ws.Range("F48").Value = 302
